$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 181, shifting rows 181:219 down to 182:220.
$ws.Rows.Item(181).Insert(-4121)

# Populate the new row 181 with the new weekly record.
$ws.Cells.Item(181, 1).Value = 4
$ws.Cells.Item(181, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(181, 3).Value = "Los Lagos"
$ws.Cells.Item(181, 4).Value = 44711
$ws.Cells.Item(181, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(181, 5).Value = 10
$ws.Cells.Item(181, 6).Value = 100112039
$ws.Cells.Item(181, 7).Value = "Ciboulette"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 80
$ws.Cells.Item(181, 11).Value = 2500
$ws.Cells.Item(181, 12).Value = 2500
$ws.Cells.Item(181, 13).Value = 2500
$ws.Cells.Item(181, 14).Value = "$/docena de atados"
$ws.Cells.Item(181, 15).Value = "Región Metropolitana"
$ws.Cells.Item(181, 16).Value = 833
$ws.Cells.Item(181, 17).Value = 3
$ws.Cells.Item(181, 18).Value = "Hortaliza"
